# Auto-generated edit script: updates cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-detected as a number by Excel;
# force Text format first so they stay stored as text (matching original inlineStr cells),
# then restore the "Normal" style so no stray formatting is left on the cell.
$numericLookingCells = @("D5", "D8", "D9", "D18", "D19", "D21", "D24", "D25", "D26", "D29", "D34", "D37", "D41", "D43", "D44", "D46", "D47", "D49")
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "28.551.91"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "1.565.85"
$ws.Range("E3").Value = "  -1.47%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "211.56"
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "46.40"
$ws.Range("D9").Value = "24.16"
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("E10").Value = "  -1.78%  "
$ws.Range("E11").Value = "  -1.34%  "
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("D13").Value = "1.790.42"
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("D14").Value = "1.574.69"
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("E15").Value = "  -2.03%  "
$ws.Range("D16").Value = "28.556.01"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("E17").Value = "  -3.11%  "
$ws.Range("D18").Value = "62.13"
$ws.Range("E18").Value = "  -1.54%  "
$ws.Range("D19").Value = "228.65"
$ws.Range("E19").Value = "  -1.69%  "
$ws.Range("E20").Value = "  -2.06%  "
$ws.Range("D21").Value = "7.33"
$ws.Range("E21").Value = "  -2.13%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  -6.08%  "
$ws.Range("D24").Value = "9.15"
$ws.Range("E24").Value = "  -2.81%  "
$ws.Range("D25").Value = "2.10"
$ws.Range("E25").Value = "  +7.24%  "
$ws.Range("D26").Value = "150.12"
$ws.Range("E26").Value = "  -1.26%  "
$ws.Range("E27").Value = "  -1.94%  "
$ws.Range("E28").Value = "  -2.65%  "
$ws.Range("D29").Value = "0.103"
$ws.Range("E29").Value = "  -3.69%  "
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("E31").Value = "  -1.76%  "
$ws.Range("E32").Value = "  -3.34%  "
$ws.Range("E33").Value = "  -1.20%  "
$ws.Range("D34").Value = "3.11"
$ws.Range("E34").Value = "  -1.57%  "
$ws.Range("D35").Value = "1.397.89"
$ws.Range("E36").Value = "  -0.94%  "
$ws.Range("D37").Value = "1.54"
$ws.Range("E37").Value = "  -2.91%  "
$ws.Range("E38").Value = "  +0.93%  "
$ws.Range("E39").Value = "  +1.97%  "
$ws.Range("E40").Value = "  -1.01%  "
$ws.Range("D41").Value = "0.537"
$ws.Range("E41").Value = "  -1.27%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "1.89"
$ws.Range("E43").Value = "  +3.05%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "0.789"
$ws.Range("E44").Value = "  -3.58%  "
$ws.Range("E45").Value = "  -4.44%  "
$ws.Range("D46").Value = "0.977"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").Value = "62.74"
$ws.Range("E47").Value = "  -2.91%  "
$ws.Range("D48").Value = "1.702.61"
$ws.Range("E48").Value = "  -1.52%  "
$ws.Range("D49").Value = "86.49"
$ws.Range("E49").Value = "  -1.28%  "
$ws.Range("E50").Value = "  -4.45%  "
$ws.Range("E51").Value = "  -1.04%  "

foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = "Normal"
}
